# Update DateBase/orders/International Ever Green_2025-10-15.xlsx
# Append new order rows (12-21) to the "Orders" sheet and refresh the
# rolled-up Number digest on the "Summary" sheet (G2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# New rows to append starting at row 12. Only columns A (PackageID), C
# (FlowerName) and F (Number) are populated, matching the existing sparse
# layout used by rows 2-11.
$newRows = @(
    @{Row=12; A="31"; C="548_白星花_tweedia white_undefined_1bunch";   F="15"},
    @{Row=13; A="";   C="572_乒乓菊白_undefined_undefined_1bunch";     F="15"},
    @{Row=14; A="";   C="734_乒乓菊红_undefined_undefined_1bunch";     F="10"},
    @{Row=15; A="32"; C="548_白星花_tweedia white_undefined_1bunch";   F="15"},
    @{Row=16; A="";   C="734_乒乓菊红_undefined_undefined_1bunch";     F="5"},
    @{Row=17; A="";   C="573_乒乓菊粉_undefined_undefined_1bunch";     F="5"},
    @{Row=18; A="";   C="418_松虫草白_scabiosa white_undefined_1bunch"; F="32"},
    @{Row=19; A="33"; C="649_洋牡丹樱花粉_undefined_undefined_1bunch"; F="20"},
    @{Row=20; A="";   C="648_洋牡丹河内_undefined_undefined_1bunch";   F="19"},
    @{Row=21; A="y";  C="";                                            F=""}
)

foreach ($item in $newRows) {
    $r = $item.Row

    if ($item.A -ne "") {
        # Force text storage so digit-only PackageIDs ("31", "32", ...)
        # don't get reinterpreted as numbers.
        $cell = $ws.Cells.Item($r, 1)
        $cell.NumberFormat = "@"
        $cell.Value = $item.A
    }

    if ($item.C -ne "") {
        $ws.Cells.Item($r, 3).Value = $item.C
    }

    if ($item.F -ne "") {
        $cell = $ws.Cells.Item($r, 6)
        $cell.NumberFormat = "@"
        $cell.Value = $item.F
    }
}

# Refresh the Summary sheet's rolled-up "Number" digest (G2) to include the
# newly appended order quantities.
$ws2 = $wb.Worksheets.Item("Summary")
$g2 = $ws2.Cells.Item(2, 7)
$g2.NumberFormat = "@"
$g2.Value = "01520520580303020515151015553220190"
